$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column (E2:E19) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E19").Value = "fullRNASEQ"

# Reflect the author's final selection in the sheet view
$ws.Range("D20:F24").Select()
